# [ADD] product_upload - agregamos el parent product
#
# Adds a new "parent" column (header in J1) to both product sheets, and
# populates the parent product code for the EINHELL sheet's single data row.

$wb = $excel.ActiveWorkbook

$wsBD = $wb.Worksheets.Item("B&D")
$wsEinhell = $wb.Worksheets.Item("EINHELL")

# New header "parent" in column J for both sheets.
$wsBD.Range("J1").Value2 = "parent"
$wsEinhell.Range("J1").Value2 = "parent"

# EINHELL's single product row gets its parent SKU populated.
$wsEinhell.Range("J2").Value2 = 4502015

# Leave the sheets with the same active-cell selections captured in the
# original edit (B&D -> J2, EINHELL -> E13), EINHELL remaining the active tab.
$wsBD.Range("J2").Select()
$wsEinhell.Activate()
$wsEinhell.Range("E13").Select()
